$d = $word.ActiveDocument

# 1. Fix typo: "willing gym battles" -> "winning gym battles"
$d.Content.Find.Execute("willing gym battles", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "winning gym battles", 2)

# 2. Add detail on how to call the elevator
$d.Content.Find.Execute("The elevator can be called (", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The elevator can be called by typing (", 2)

# 3. Remove the stray _GoBack bookmark from its old location
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 4. Fix capitalization typo: "Eventually, Pogomen" -> "Eventually, pogomen"
$d.Content.Find.Execute("Eventually, Pogomen", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Eventually, pogomen", 2)

# 5. Expand on the ending text, then move the _GoBack bookmark to mark
#    this as the most-recently-edited location.
$target = $d.Content
$target.Find.Execute("check out each ending.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, `
                      "check out other endings (although some depend on choices made earlier in the story and may not be reachable from this point).", 2)

$endRange = $d.Content
$endRange.Find.Execute("may not be reachable from this point).", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)
$bmRange = $d.Range($endRange.End, $endRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
